# Add support for phonetic pronunciation
# Adds a new "phoneme pronunciation" sheet with an example row (ThioJoe -> IPA)

$wb = $excel.ActiveWorkbook

# Add the new worksheet at the end of the workbook (After the last existing sheet)
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "phoneme pronunciation"

# Header row
$newSheet.Range("A1").Value = "Text"
$newSheet.Range("B1").Value = "Phonetic Pronunciation"
$newSheet.Range("C1").Value = "Case Sensitive (True/False)"
$newSheet.Range("D1").Value = "Phonetic Alphabet"

# Example data row
$newSheet.Range("A2").Value = "ThioJoe"
$newSheet.Range("B2").Value = "ˈθioʊd͡ʒoʊ"
$newSheet.Range("C2").Value = $false
$newSheet.Range("D2").Value = "ipa"

# Make the new sheet the active/selected tab
$newSheet.Activate()
$newSheet.Range("D6").Select()
